$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.948.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.508.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.53"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.16"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.92%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.41"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.954.33"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "58.849.83"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.72"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.34%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.523.87"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "322.26"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.07"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.44%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.163"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.47%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.55"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0762"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.48"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.55%  "

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.26"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.42%  "

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.75"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.35%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.02%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.36"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.03"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.28%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.57"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.799"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "281.05"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.56%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.38%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.603"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.33%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "129.75"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.54%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.99"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.90"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.22%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0500"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.66%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.22"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.761.12"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.982"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.60%  "

Write-Host "Edit complete"
